$wb = $excel.ActiveWorkbook

# --- Sheet1 (Log): add rows 47 and 48 ---
$wsLog = $wb.Worksheets.Item("Log")

$wsLog.Range("A47").Value = "27/02/2026"
$wsLog.Range("B47").Value = "14:00"
$wsLog.Range("C47").Value = "Proyección 3 meses e Int. por caución proyectado"
$wsLog.Range("D47").Value = "Configuración (Configuración en menú): método Mediana/Promedio y meses de historia (3, 6, 12, 24). Próximos 3 meses proyectados con ventana rodante. Int. por caución: punto de partida = último mes real (G/P + interés), luego última tasa conocida aplicada en cadena para cada mes proyectado."
$wsLog.Range("E47").Value = "Implementacion"

$wsLog.Range("A48").Value = "27/02/2026"
$wsLog.Range("B48").Value = "14:15"
$wsLog.Range("C48").Value = "Disclaimer bajo proyección"
$wsLog.Range("D48").Value = "Texto en letra chica y gris oscuro bajo la proyección indicando metodología: Mediana/Promedio de N meses, ventana rodante, y cómo se calcula Int. por caución proyectado."
$wsLog.Range("E48").Value = "Implementacion"

# --- Sheet2 (Resumen): add rows 43, 44, 45 ---
$wsResumen = $wb.Worksheets.Item("Resumen")

$wsResumen.Range("A43").Value = "Proyección próximos 3 meses"
$wsResumen.Range("B43").Value = "Debajo del total real en Flujo por mes: ""Próximos 3 meses proyectados"" con ventana rodante. Configuración (sidebar): método (Mediana/Promedio) y meses de historia (3, 6, 12, 24). Ingresos, egresos, G/P y ratios proyectados por mes."

$wsResumen.Range("A44").Value = "Int. por caución proyectado"
$wsResumen.Range("B44").Value = "Para cada mes proyectado: punto de partida = G/P + interés del período anterior (último real para mes 1; proyectado 1 para mes 2; proyectado 2 para mes 3). Se aplica la última tasa conocida en cadena. Sin salto respecto al último valor real."

$wsResumen.Range("A45").Value = "Disclaimer proyección"
$wsResumen.Range("B45").Value = "Debajo de las filas proyectadas, texto en letra chica y gris oscuro que explica la metodología: Mediana/Promedio de N meses, ventana rodante, y cálculo de Int. por caución proyectado."

# --- Sheet4 (Versiones): add row 18 ---
$wsVersiones = $wb.Worksheets.Item("Versiones")

# Force A18 to be stored as text (matches existing version numbers like "1.15"
# which are text, not numeric) by setting the cell format to Text before
# assigning the value.
$wsVersiones.Range("A18").NumberFormat = "@"
$wsVersiones.Range("A18").Value = "1.16"
$wsVersiones.Range("B18").Value = "27/02/2026"
$wsVersiones.Range("C18").Value = "Proyección 3 meses: config (mediana/promedio, meses historia), ventana rodante; Int. por caución proyectado con punto de partida = último real (G/P+interés) y última tasa en cadena; disclaimer bajo proyección. Despliegue a producción."
